$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("测试需求分析表")

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "完成对前端页面响应式布局的测试内容"
$ws.Range("E6").Value = "无"
$ws.Range("F6").Value = "王康明"

$ws.Rows.Item(6).RowHeight = $ws.Rows.Item(5).RowHeight

$ws.Select()
$ws.Range("D18").Select()
